$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Simple header-field value updates
# ---------------------------------------------------------------------------
# CO No. changed
$ws.Range("C14").Value = "cs1709-2310001"

# Order Date changed (18 Oct 2023 -> 26 Oct 2023, serial 45217 -> 45225)
$ws.Range("C18").Value = 45225

# Target Date text changed
$ws.Range("C19").Value = "01 Nov 2023 - 10 Nov 2023"

# Forecast period label changed
$ws.Range("N22").Value = "11 Nov ~ 20 Nov"

# Forecast/Inbound-plan date changed (18 Oct 2023 -> 26 Oct 2023)
$ws.Range("P23").Value = 45225

# ---------------------------------------------------------------------------
# 2. Reorder the order-line detail grid (rows 24-28).
#    Mapping of new row <- old row content:
#      24 <- 28
#      25 <- 24
#      26 <- 26 (unchanged)
#      27 <- 27 (unchanged)
#      28 <- 25
# ---------------------------------------------------------------------------
$gridCols = @("B","C","D","F","G","H","I","J","K","L","M","N","O","Q","R")

function Get-RowVals {
    param($ws, $r, $cols)
    $result = @{}
    foreach ($c in $cols) {
        $result[$c] = $ws.Range("$c$r").Value2
    }
    return $result
}

function Set-RowVals {
    param($ws, $r, $cols, $vals)
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $vals[$c]
    }
}

# Snapshot the current (pre-edit) content of the three rows involved in the cycle
$old24 = Get-RowVals $ws 24 $gridCols
$old25 = Get-RowVals $ws 25 $gridCols
$old28 = Get-RowVals $ws 28 $gridCols

# Write back in the new order
Set-RowVals $ws 24 $gridCols $old28
Set-RowVals $ws 25 $gridCols $old24
Set-RowVals $ws 28 $gridCols $old25
